$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attributes")

# Add two new attribute rows for the "rdcomponents_institutions" entity,
# growing the example dataset:
#   province -> string
#   type     -> string
$ws.Range("A16").Value = "rdcomponents_institutions"
$ws.Range("B16").Value = "province"
$ws.Range("C16").Value = "string"

$ws.Range("A17").Value = "rdcomponents_institutions"
$ws.Range("B17").Value = "type"
$ws.Range("C17").Value = "string"

# Match the existing "name" column formatting (explicit black font) used
# by every other row in this table.
$ws.Range("B16").Font.Color = 0
$ws.Range("B17").Font.Color = 0

$ws.Range("B17").Select()
